$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 282, shifting existing rows 282-376 down to 283-377.
# Excel extends the formatting of the row above into the newly inserted row,
# which matches the style (s="2") already present on column D of the
# surrounding rows.
$ws.Rows.Item(282).Insert()

# Populate the carried-over (unchanged) columns of the new row 282 by copying
# them from the row that is now directly below it (the old row 282, now row 283).
$ws.Range("A282").Value = $ws.Range("A283").Value2
$ws.Range("B282").Value = $ws.Range("B283").Value2
$ws.Range("C282").Value = $ws.Range("C283").Value2
$ws.Range("E282").Value = $ws.Range("E283").Value2
$ws.Range("F282").Value = $ws.Range("F283").Value2
$ws.Range("G282").Value = $ws.Range("G283").Value2
$ws.Range("H282").Value = $ws.Range("H283").Value2
$ws.Range("I282").Value = $ws.Range("I283").Value2
$ws.Range("N282").Value = $ws.Range("N283").Value2
$ws.Range("O282").Value = $ws.Range("O283").Value2
$ws.Range("Q282").Value = $ws.Range("Q283").Value2
$ws.Range("R282").Value = $ws.Range("R283").Value2

# Set the new data values for the inserted row 282.
$ws.Range("D282").Value = 45229
$ws.Range("J282").Value = 230
$ws.Range("K282").Value = 9000
$ws.Range("L282").Value = 10000
$ws.Range("M282").Value = 9435
$ws.Range("P282").Value = 189
